$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 71431576
$ws.Range("I18").Value = 2989
$ws.Range("K18").Value = 2989
$ws.Range("M18").Value = -2705

$ws.Range("H55").Value = 822.3333
$ws.Range("J55").Value = 808
$ws.Range("L55").Value = 808
$ws.Range("N55").Value = -1236

$ws.Range("H70").Value = 2973.5
$ws.Range("I70").Value = 3002
$ws.Range("J70").Value = 2964
$ws.Range("K70").Value = 9006
$ws.Range("L70").Value = 8892
$ws.Range("M70").Value = -8736
$ws.Range("N70").Value = -9432

$ws.Range("H73").Value = 2973.5
$ws.Range("I73").Value = 3002
$ws.Range("J73").Value = 2964
$ws.Range("K73").Value = 9006
$ws.Range("L73").Value = 8892
$ws.Range("M73").Value = -8070
$ws.Range("N73").Value = -10764

$ws.Range("H96").Value = 401.0909
$ws.Range("I96").Value = 334.77777
$ws.Range("K96").Value = 1004.33331
$ws.Range("M96").Value = 368.66669

$ws.Range("H116").Value = 4811
$ws.Range("J116").Value = 8999.25
$ws.Range("L116").Value = 8999.25
$ws.Range("N116").Value = -15883.25


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3079.6667
$ws.Range("I32").Value = 2880.9077
$ws.Range("K32").Value = 2880.9077
$ws.Range("M32").Value = -2593.9077

$ws.Range("H61").Value = 3113
$ws.Range("I61").Value = 1810.2858
$ws.Range("K61").Value = 1810.2858
$ws.Range("M61").Value = -1598.2858

$ws.Range("H74").Value = 102133.49
$ws.Range("I74").Value = 130088.07
$ws.Range("K74").Value = 130088.07
$ws.Range("M74").Value = -129214.07

$ws.Range("H77").Value = 102133.49
$ws.Range("I77").Value = 130088.07
$ws.Range("K77").Value = 650440.3500000001
$ws.Range("M77").Value = -646072.3500000001

$ws.Range("H110").Value = 2797.2368
$ws.Range("I110").Value = 1858
$ws.Range("J110").Value = 8996.200000000001
$ws.Range("K110").Value = 1858
$ws.Range("L110").Value = 8996.200000000001
$ws.Range("M110").Value = 187
$ws.Range("N110").Value = -13086.2

$ws.Range("H136").Value = 3113
$ws.Range("I136").Value = 1810.2858
$ws.Range("K136").Value = 5430.857400000001
$ws.Range("M136").Value = -2880.857400000001


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4784.5
$ws.Range("I86").Value = 3755.8
$ws.Range("K86").Value = 3755.8
$ws.Range("M86").Value = -2632.8

$ws.Range("H89").Value = 4784.5
$ws.Range("I89").Value = 3755.8
$ws.Range("K89").Value = 18779
$ws.Range("M89").Value = -13163

$ws.Range("H105").Value = 23639646
$ws.Range("I105").Value = 1669402.4
$ws.Range("K105").Value = 1669402.4
$ws.Range("M105").Value = -1667655.4

$ws.Range("H107").Value = 3078624.5
$ws.Range("I107").Value = 4275041.5
$ws.Range("K107").Value = 4275041.5
$ws.Range("M107").Value = -4273121.5

$ws.Range("H134").Value = 2850.2
$ws.Range("I134").Value = 2667.3333
$ws.Range("J134").Value = 2928.5715
$ws.Range("K134").Value = 8001.999899999999
$ws.Range("L134").Value = 8785.7145
$ws.Range("M134").Value = -5466.999899999999
$ws.Range("N134").Value = -13855.7145


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1593.6666
$ws.Range("J16").Value = 1728.0834
$ws.Range("L16").Value = 1728.0834
$ws.Range("N16").Value = -2302.0834

$ws.Range("H58").Value = 2482.7297
$ws.Range("I58").Value = 1696.76
$ws.Range("K58").Value = 1696.76
$ws.Range("M58").Value = -1493.76

$ws.Range("H113").Value = 1593.6666
$ws.Range("J113").Value = 1728.0834
$ws.Range("L113").Value = 1728.0834
$ws.Range("N113").Value = -6068.0834

$ws.Range("H132").Value = 2964.4138
$ws.Range("I132").Value = 2225.9524
$ws.Range("J132").Value = 4902.875
$ws.Range("K132").Value = 6677.8572
$ws.Range("L132").Value = 14708.625
$ws.Range("M132").Value = -4147.8572
$ws.Range("N132").Value = -19768.625

$ws.Range("H136").Value = 2482.7297
$ws.Range("I136").Value = 1696.76
$ws.Range("K136").Value = 5090.28
$ws.Range("M136").Value = -2540.28


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 167221.33
$ws.Range("J121").Value = 334366.34
$ws.Range("L121").Value = 1003099.02
$ws.Range("N121").Value = -1005719.02

$ws.Range("H129").Value = 964.5
$ws.Range("J129").Value = 2388
$ws.Range("L129").Value = 7164
$ws.Range("N129").Value = -17164

$ws.Range("H131").Value = 7579.32
$ws.Range("J131").Value = 2376.6155
$ws.Range("L131").Value = 7129.8465
$ws.Range("N131").Value = -17209.8465

$ws.Range("H134").Value = 15389152
$ws.Range("I134").Value = 18183542
$ws.Range("K134").Value = 54550626
$ws.Range("M134").Value = -54545556


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3785.8928
$ws.Range("I132").Value = 1192.5385
$ws.Range("K132").Value = 3577.6155
$ws.Range("M132").Value = -1047.6155


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 44357
$ws.Range("I40").Value = 49017.13
$ws.Range("K40").Value = 49017.13
$ws.Range("M40").Value = -48881.13

$ws.Range("H128").Value = 0
$ws.Range("J128").Value = 0
$ws.Range("L128").Value = 0
$ws.Range("N128").ClearContents()

$ws.Range("H132").Value = 6248.909
$ws.Range("J132").Value = 9429.143
$ws.Range("L132").Value = 28287.429
$ws.Range("N132").Value = -33347.429

$ws.Range("H140").Value = 77596
$ws.Range("J140").Value = 77596
$ws.Range("L140").Value = 77596
$ws.Range("N140").Value = -87956


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 64695.855
$ws.Range("J46").Value = 64695.855
$ws.Range("L46").Value = 64695.855
$ws.Range("N46").Value = -65157.855

$ws.Range("H113").Value = 585.1
$ws.Range("I113").Value = 543.75
$ws.Range("J113").Value = 647.125
$ws.Range("K113").Value = 1631.25
$ws.Range("L113").Value = 1941.375
$ws.Range("M113").Value = 538.75
$ws.Range("N113").Value = -6281.375

$ws.Range("H122").Value = 8067881.5
$ws.Range("J122").Value = 31252136
$ws.Range("L122").Value = 93756408
$ws.Range("N122").Value = -93761308

$ws.Range("H134").Value = 64695.855
$ws.Range("J134").Value = 64695.855
$ws.Range("L134").Value = 194087.565
$ws.Range("N134").Value = -199157.565

$ws.Range("H136").Value = 37040050
$ws.Range("I136").Value = 45455732
$ws.Range("K136").Value = 136367196
$ws.Range("M136").Value = -136364646
